$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 (b.md.md) status changes from
#     "Handed back: in sync with en-US" to "Ready for handoff"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: row 3 (b.md.md) gets a fresh handoff
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
$wsZh.Range("D3").Value = "2016-01-28 10:57:02"

# Update the existing hyperlink's display text in place (keeps the same
# relationship/target, only the visible text changes) - the Hyperlinks
# collection's indexer is unreliable in this host, so walk it instead.
$i = 0
foreach ($h in $wsZh.Hyperlinks) {
    $i++
    if ($i -eq 6) {
        $h.TextToDisplay = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
    }
}

# --- de-de sheet: row 3 (b.md.md) gets a fresh handoff
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
$wsDe.Range("D3").Value = "2016-01-28 10:57:15"

$j = 0
foreach ($h in $wsDe.Hyperlinks) {
    $j++
    if ($j -eq 6) {
        $h.TextToDisplay = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
    }
}
